$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "Computes Hotelling's multivariate t-test for each partially observed covariate, examining patient differences conditional on having an observed covariate value or not."

$ws.Range("D3").Select()
